$d = $word.ActiveDocument

# Locate the target question paragraph ("What is your favorite mobile app
# platform to develop in? Explain.") and the empty paragraph immediately
# preceding it.
$target = $null
$prev = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "What is your favorite mobile app platform*") {
        $target = $p
        break
    }
    $prev = $p
}

# Move the `_GoBack` bookmark onto the preceding (empty) paragraph. Adding a
# bookmark with a name that already exists elsewhere in the document moves
# it, removing the old occurrence automatically.
if ($prev -ne $null) {
    $d.Bookmarks.Add("_GoBack", $prev.Range)
}

# Replace the whole (multi-run) question text with the new single run of
# text.
$target.Range.Find.Execute(
    "What is your favorite mobile app platform to develop in? Explain.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Would you like to see more mobile application development classes offered at Rose?",
    2)
